$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": M21 1244.13 -> 1632.93
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("M21").Value = 1632.93

# Sheet "VENTA MENSUAL": F21 1244.13 -> 1632.93, F26 20529.03 -> 20917.83
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F21").Value = 1632.93
$wsVentaMensual.Range("F26").Value = 20917.83

# Sheet "CUMPLIMIENTO MENSUAL": D12/E12/F12 and D15/E15/F15 updated
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 20665.91
$wsCumplimiento.Range("E12").Value = 22434.1754117774
$wsCumplimiento.Range("F12").Value = 0.4794865207935968

$wsCumplimiento.Range("D15").Value = 20917.83
$wsCumplimiento.Range("E15").Value = 37285.63623249458
$wsCumplimiento.Range("F15").Value = 0.3593914822262205
